$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API's V2")

# --- Populate new cell values in the exact order the strings were first
# --- introduced in the target workbook, so that sharedStrings.xml ends up
# --- with the same append order as the authored edit.

# Row 9 - create group
$ws.Range("A9").Value = "http://localhost:8080/group/create"
$ws.Range("D9").Value = @'
{
    "groupName": "NewGroup112",
    "isPublic": true,
    "createdBy": 1,
    "createdDate": null,
    "isActive": null
}
'@
$ws.Range("E9").Value = @'
{
    "message": "Created Successfully",
    "results": {
        "groupId": 52,
        "groupName": "NewGroup112",
        "isPublic": true,
        "createdBy": 1,
        "createdDate": "2020-11-28T12:31:27.000+00:00",
        "isActive": null
    }
}
'@

# Row 10 - update group
$ws.Range("A10").Value = "http://localhost:8080/group/update"
$ws.Range("D10").Value = @'
{
    "groupId": 51,
    "groupName": "NewGroup1132",
    "isPublic": false,
    "createdBy": 1,
    "createdDate": null,
    "isActive": null
}
'@
$ws.Range("E10").Value = @'
{
    "message": "Updated Successfully",
    "results": {
        "groupId": 51,
        "groupName": "NewGroup1132",
        "isPublic": false,
        "createdBy": 1,
        "createdDate": "2020-11-28T12:10:44.000+00:00",
        "isActive": null
    }
}
'@

# Row 11 - delete group
$ws.Range("A11").Value = "http://localhost:8080/group/delete"
$ws.Range("C11").Value = "To Delete existing group/Board"
$ws.Range("C10").Value = "To updae existing group/Board"
$ws.Range("C9").Value = "To add a new group/Board"
$ws.Range("D11").Value = @'
{
    "groupId": 51
}
'@

# Row 12 - get owner groups
$ws.Range("A12").Value = "http://localhost:8080/group/getOwnerGroups"
$ws.Range("C12").Value = "To Get all Board owner groups "
$ws.Range("E12").Value = @'
[    {
        "groupId": 48,
        "groupName": "Test-Board",
        "isPublic": true,
        "createdBy": 1,
        "createdDate": null,
        "isActive": null
    },
    {
        "groupId": 50,
        "groupName": "Date Check",
        "isPublic": true,
        "createdBy": 1,
        "createdDate": null,
        "isActive": null
    },
    {
        "groupId": 51,
        "groupName": "NewGroup1132",
        "isPublic": false,
        "createdBy": 1,
        "createdDate": "2020-11-28T12:10:44.000+00:00",
        "isActive": false
    },
    {
        "groupId": 52,
        "groupName": "NewGroup112",
        "isPublic": true,
        "createdBy": 1,
        "createdDate": "2020-11-28T12:31:27.000+00:00",
        "isActive": null
    }
]
'@

# Remaining columns reusing already-known shared strings.
$ws.Range("B9").Value = "POST"
$ws.Range("B10").Value = "POST"
$ws.Range("B11").Value = "POST"
$ws.Range("B12").Value = "POST"
$ws.Range("E11").Value = $ws.Range("E7").Value2
$ws.Range("D12").Value = $ws.Range("D8").Value2

# --- Formatting: row heights ---
$ws.Rows.Item(9).RowHeight = 165
$ws.Rows.Item(10).RowHeight = 165
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 409.5

# --- Hyperlinks (must be added before final style tweaks so the
# --- "Hyperlink" cell style used below resolves to the already-existing
# --- style slots instead of leaving stray duplicates behind) ---
$ws.Hyperlinks.Add($ws.Range("A9"), "http://localhost:8080/group/create")
$ws.Hyperlinks.Add($ws.Range("A10"), "http://localhost:8080/group/update")
$ws.Hyperlinks.Add($ws.Range("A11"), "http://localhost:8080/group/delete")
$ws.Hyperlinks.Add($ws.Range("A12"), "http://localhost:8080/group/getOwnerGroups")

# --- Formatting: cell styles ---
$ws.Range("A9").VerticalAlignment = -4108
$ws.Range("A10").VerticalAlignment = -4108
$ws.Range("A11").Style = "Hyperlink"
$ws.Range("A12").VerticalAlignment = -4108

$ws.Range("B9:C10").VerticalAlignment = -4108
$ws.Range("B11:C12").VerticalAlignment = -4108

$ws.Range("D9:E10").VerticalAlignment = -4108
$ws.Range("D9:E10").WrapText = $true
$ws.Range("D11:E11").VerticalAlignment = -4108
$ws.Range("D11:E11").WrapText = $true
$ws.Range("D12").VerticalAlignment = -4108
$ws.Range("D12").WrapText = $true
$ws.Range("E12").WrapText = $true

# --- Final selection matches the authored edit ---
$ws.Range("E12").Select()
